$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 375 (pushes old rows 375..471 down to 376..472).
# Excel copies the formatting of the row above the insertion point (row 374),
# which already gives column D the correct date style (s="2").
$ws.Rows.Item(375).Insert()

# Populate the newly inserted row 375 with a new price entry for Lechuga /
# Conconina(o) - Primera, duplicating the existing 200/4000/4500/4250/425
# record (previously at row 383, now shifted to row 384) but dated 44551.
$ws.Cells.Item(375, 1).Value = 7
$ws.Cells.Item(375, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(375, 3).Value = "Ñuble"
$ws.Cells.Item(375, 4).Value = 44551
$ws.Cells.Item(375, 5).Value = 16
$ws.Cells.Item(375, 6).Value = 100112033
$ws.Cells.Item(375, 7).Value = "Lechuga"
$ws.Cells.Item(375, 8).Value = "Conconina(o)"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 200
$ws.Cells.Item(375, 11).Value = 4000
$ws.Cells.Item(375, 12).Value = 4500
$ws.Cells.Item(375, 13).Value = 4250
$ws.Cells.Item(375, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(375, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(375, 16).Value = 425
$ws.Cells.Item(375, 17).Value = 10
$ws.Cells.Item(375, 18).Value = "Hortaliza"
